$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-row price / volume(1h) updates (column D values are forced to Text
# format first since the source sheet stores them as inline strings, e.g.
# "1.002", "13.00", "0.00001229" - values Excel would otherwise coerce to
# numbers and reformat).
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '23.706.98'
$ws.Range('E2').Value = '  +0.85%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.657.62'
$ws.Range('E3').Value = '  +1.20%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('E5').Value = '  -0.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '303.01'
$ws.Range('E6').Value = '  -0.26%  '
$ws.Range('E7').Value = '  +0.77%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3613'
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '51.28'
$ws.Range('E9').Value = '  -0.67%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.08191'
$ws.Range('E10').Value = '  +0.36%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.230'
$ws.Range('E11').Value = '  -0.05%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.002'
$ws.Range('E12').Value = '  +0.11%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.53'
$ws.Range('E13').Value = '  +0.32%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.473'
$ws.Range('E14').Value = '  +0.33%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.406'
$ws.Range('E15').Value = '  +0.59%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001229'
$ws.Range('E16').Value = '  -0.75%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.652.89'
$ws.Range('E17').Value = '  +0.90%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '97.73'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.07009'
$ws.Range('E19').Value = '  +1.00%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.823'
$ws.Range('E20').Value = '  +3.63%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.64'
$ws.Range('E21').Value = '  +0.73%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.000'
$ws.Range('E22').Value = '  +0.00%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.80'
$ws.Range('E23').Value = '  +2.49%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '23.714.76'
$ws.Range('E24').Value = '  +0.93%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.506'
$ws.Range('E25').Value = '  +0.64%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.013'
$ws.Range('E26').Value = '  -1.33%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '21.25'
$ws.Range('E27').Value = '  +0.47%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '153.81'
$ws.Range('E28').Value = '  +1.37%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.225'
$ws.Range('E29').Value = '  -0.49%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '134.24'
$ws.Range('E30').Value = '  +1.08%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.838.62'
$ws.Range('E31').Value = '  +1.46%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.044'
$ws.Range('E32').Value = '  +6.97%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.227'
$ws.Range('E33').Value = '  +2.92%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '12.04'
$ws.Range('E34').Value = '  +5.42%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.059'
$ws.Range('E35').Value = '  -3.37%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02806'
$ws.Range('E36').Value = '  +1.69%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2515'
$ws.Range('E37').Value = '  +0.72%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.08798'
$ws.Range('E38').Value = '  +0.36%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.083'
$ws.Range('E39').Value = '  +1.66%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.07029'
$ws.Range('E40').Value = '  -0.71%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '13.00'
$ws.Range('E41').Value = '  +6.64%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6993'
$ws.Range('E42').Value = '  -0.65%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.336'
$ws.Range('E43').Value = '  -1.10%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.14'
$ws.Range('E44').Value = '  +3.79%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6517'
$ws.Range('E45').Value = '  -0.11%  '
$ws.Range('E48').Value = '  -0.01%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.07929'
$ws.Range('E49').Value = '  -0.51%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '127.99'
$ws.Range('E50').Value = '  -0.86%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.183'
$ws.Range('E51').Value = '  -0.67%  '

# Rows 46/47 swap (NEARProtocol <-> Frax) plus updated price/volume
$ws.Range('B46').Value = 'Frax'
$ws.Range('C46').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.9995'
$ws.Range('E46').Value = '  -0.08%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.309'
$ws.Range('E47').Value = '  +1.37%  '
